# Add a new "Complaint - creator read access" rule row to Sheet1, mirroring
# the existing "Case File - creator read access" row (row 26) but scoped to
# the COMPLAINT object type instead of CASE_FILE.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Duplicate row 26 (styles + the "reader, creator" action in column H) into
# the new row 27.
$ws.Range("A26:H26").Copy($ws.Range("A27:H27"))

# Row 26 is a tall, word-wrapped row (auto height); give the new row the
# same visible height.
$ws.Rows.Item(27).RowHeight = $ws.Rows.Item(26).RowHeight()

# Differentiate the new rule from the one it was copied from.
$ws.Range("B27").Value = "Complaint - creator read access"
$ws.Range("C27").Value = "COMPLAINT"
